# Add dictionary configuration functionality:
# Create a new "lkWordDict" worksheet after the existing "lkMetaCategories"
# sheet, with header columns "Translation" (B1) and "PoS" (C1).

$wb = $excel.ActiveWorkbook

$metaSheet = $wb.Worksheets.Item("lkMetaCategories")

# Insert the new sheet right after lkMetaCategories so ordering + the
# workbook's activeTab / tabSelected bookkeeping line up with the edit.
$wordDict = $wb.Worksheets.Add($null, $metaSheet)
$wordDict.Name = "lkWordDict"

# Header row
$wordDict.Range("B1").Value = "Translation"
$wordDict.Range("C1").Value = "PoS"

# Size column B to fit its header text, and leave the active cell on A2
# (matching the selection saved with the sheet).
$wordDict.Columns.Item(2).ColumnWidth = 8.7
$wordDict.Range("A2").Select() | Out-Null

$wordDict.Activate()
